$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 196, shifting existing rows 196:249 down to 197:250.
$ws.Rows.Item(196).Insert()

# Populate the newly inserted row 196 with the new record's data.
$ws.Cells.Item(196, 1).Value = 5
$ws.Cells.Item(196, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(196, 3).Value = "Maule"
$ws.Cells.Item(196, 4).Value = 44642
$ws.Cells.Item(196, 5).Value = 7
$ws.Cells.Item(196, 6).Value = "Fruta"
$ws.Cells.Item(196, 7).Value = 100102
$ws.Cells.Item(196, 8).Value = "Cítricos"
$ws.Cells.Item(196, 9).Value = 100102004
$ws.Cells.Item(196, 10).Value = "Mandarina"
$ws.Cells.Item(196, 11).Value = "Murcott"
$ws.Cells.Item(196, 12).Value = "Primera"
$ws.Cells.Item(196, 13).Value = 210
$ws.Cells.Item(196, 14).Value = 11000
$ws.Cells.Item(196, 15).Value = 11000
$ws.Cells.Item(196, 16).Value = 11000
$ws.Cells.Item(196, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(196, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(196, 19).Value = 611
$ws.Cells.Item(196, 20).Value = 18
